# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Updates the "Periodo Mora" table (rows 16-50) on Hoja1:
#  - Inserts a new worker (ESTIBENSON RAFAEL PALACIN VILLAREAL) with period 2211
#  - Re-sequences the existing recurring workers (ERICH RAFAEL HERRERA CABALLERO,
#    GUILLERMO RODRIGUEZ PIÑERES, CARLOSMARIO SEPULVEDA PEDROZO, LEIDY MARIAM
#    BELLIDO ARROYO) across periods 2302-2401, updating "Valor Mora" amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora,
#          F=Valor Mora, G=Salario Basico
$rows = @(
    @(16, "CC", "1193522083", "ESTIBENSON RAFAEL PALACIN VILLAREAL", "2211", 28000,  1000000),
    @(17, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2302", 34027,  1160000),
    @(18, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2302", 15467,  1160000),
    @(19, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2303", 46400,  1160000),
    @(20, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2303", 46400,  1160000),
    @(21, "CC", "1143412227", "LEIDY MARIAM BELLIDO ARROYO",         "2303", 1547,   1160000),
    @(22, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2304", 46400,  1160000),
    @(23, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2304", 46400,  1160000),
    @(24, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2305", 46400,  1160000),
    @(25, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2305", 38667,  1160000),
    @(26, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2305", 46400,  1160000),
    @(27, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2306", 46400,  1160000),
    @(28, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2306", 46400,  1160000),
    @(29, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2306", 46400,  1160000),
    @(30, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2307", 46400,  1160000),
    @(31, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2307", 46400,  1160000),
    @(32, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2307", 46400,  1160000),
    @(33, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2308", 46400,  1160000),
    @(34, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2308", 46400,  1160000),
    @(35, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2308", 46400,  1160000),
    @(36, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2309", 46400,  1160000),
    @(37, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2309", 46400,  1160000),
    @(38, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2309", 46400,  1160000),
    @(39, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2310", 46400,  1160000),
    @(40, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2310", 46400,  1160000),
    @(41, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2310", 46400,  1160000),
    @(42, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2311", 46400,  1160000),
    @(43, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2311", 46400,  1160000),
    @(44, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2311", 46400,  1160000),
    @(45, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2312", 46400,  1160000),
    @(46, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2312", 46400,  1160000),
    @(47, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2312", 46400,  1160000),
    @(48, "CC", "9315662",    "ERICH RAFAEL HERRERA CABALLERO",      "2401", 46400,  1160000),
    @(49, "CC", "1052219882", "CARLOSMARIO SEPULVEDA PEDROZO",       "2401", 46400,  1160000),
    @(50, "CC", "73119562",   "GUILLERMO RODRIGUEZ PIÑERES",         "2401", 46400,  1160000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($rowNum, 3).Value = $r[2]   # C: N Doc Trabajador
    $ws.Cells.Item($rowNum, 4).Value = $r[3]   # D: Nombre Trabajador
    $ws.Cells.Item($rowNum, 5).Value = $r[4]   # E: Periodo Mora
    $ws.Cells.Item($rowNum, 6).Value = $r[5]   # F: Valor Mora
    $ws.Cells.Item($rowNum, 7).Value = $r[6]   # G: Salario Basico
}
